$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new column (K) for year 2023, copying the
# formatting from the preceding column (J) for each of the existing rows,
# then filling in the new data values.

# Row 3: year header
$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)
$ws.Range("K3").Value = 2023

# Row 4: total average monthly remuneration
$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial(-4122)
$ws.Range("K4").Value = 602.20000000000005

# Row 5: women
$ws.Range("J5").Copy()
$ws.Range("K5").PasteSpecial(-4122)
$ws.Range("K5").Value = 337.9

# Row 6: men
$ws.Range("J6").Copy()
$ws.Range("K6").PasteSpecial(-4122)
$ws.Range("K6").Value = 666.1

$excel.CutCopyMode = 0
